$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "스크립트"

# --- Update cell values ---
# Row 1: title label + new title value
$ws.Range("A1").Value = "제목"
$ws.Range("B1").Value = "강사 메일 스크립트 자동 생성 Input Parameters"

# Row 2: description label (replaces old placeholder text), value left blank
$ws.Range("A2").Value = "설명"
$ws.Range("B2").Value = ""

# Row 3: blank spacer row
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = ""

# Row 4: unchanged values, just restyled further below
$ws.Range("A4").Value = "부문"
$ws.Range("B4").Value = "구매·자재"

# Row 5 / Row 6 values stay the same as before (greeting / closing text)

# Row 7: new row for the source schedule file name
$ws.Range("A7").Value = "일정표 원본파일명"
$ws.Range("B7").Value = "2019일정계획표(2019.08.02).xlsx"

# --- Column width for column A ---
$ws.Columns.Item(1).ColumnWidth = 16.2857142857143

# --- Apply thin box border + no alignment to A1:B3 and B7 ---
$rNoAlign1 = $ws.Range("A1:B3")
$rNoAlign1.Borders.LineStyle = 1
$rNoAlign1.Borders.Weight = 2

$rNoAlign2 = $ws.Range("B7")
$rNoAlign2.Borders.LineStyle = 1
$rNoAlign2.Borders.Weight = 2

# --- Apply thin box border + vertical-center alignment to A4:B4, A5, A6, A7 ---
$rValign1 = $ws.Range("A4:B4")
$rValign1.Borders.LineStyle = 1
$rValign1.Borders.Weight = 2
$rValign1.VerticalAlignment = -4108

$rValign2 = $ws.Range("A5")
$rValign2.Borders.LineStyle = 1
$rValign2.Borders.Weight = 2
$rValign2.VerticalAlignment = -4108

$rValign3 = $ws.Range("A6")
$rValign3.Borders.LineStyle = 1
$rValign3.Borders.Weight = 2
$rValign3.VerticalAlignment = -4108

$rValign4 = $ws.Range("A7")
$rValign4.Borders.LineStyle = 1
$rValign4.Borders.Weight = 2
$rValign4.VerticalAlignment = -4108

# --- Apply thin box border + vertical-center + wrap to B5:B6 ---
$rWrap = $ws.Range("B5:B6")
$rWrap.Borders.LineStyle = 1
$rWrap.Borders.Weight = 2
$rWrap.VerticalAlignment = -4108
$rWrap.WrapText = $true

# --- Selection matches the final saved cursor position ---
$ws.Range("B13").Select()
